$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,14
$data[0,0] = 0.4047206666666667
$data[0,1] = 1.214162
$data[0,2] = 0.02578034976888792
$data[0,3] = 0.02578034976888792
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.01779166666666667
$data[0,7] = 0.053375
$data[0,8] = 0.03510503888026929
$data[0,9] = 0.03510503888026929
$data[0,10] = 0.007200655194444445
$data[0,11] = 0.06480589674999999
$data[0,12] = 0.0009050201809837519
$data[0,13] = 0.0009050201809837517
$data[1,0] = 0.4047206666666667
$data[1,1] = 1.214162
$data[1,2] = 0.02578034976888792
$data[1,3] = 0.02578034976888792
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.036329
$data[1,7] = 0.108987
$data[1,8] = 0.07168136529168917
$data[1,9] = 0.07168136529168917
$data[1,10] = 0.01470309709933333
$data[1,11] = 0.132327873894
$data[1,12] = 0.00184797066913117
$data[1,13] = 0.00184797066913117
$data[2,0] = 0.4047206666666667
$data[2,1] = 1.214162
$data[2,2] = 0.02578034976888792
$data[2,3] = 0.02578034976888792
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.09574300000000001
$data[2,7] = 0.287229
$data[2,8] = 0.1889121351295713
$data[2,9] = 0.1889121351295713
$data[2,10] = 0.03874917078866667
$data[2,11] = 0.348742537098
$data[2,12] = 0.004870220919227767
$data[2,13] = 0.004870220919227767
$data[3,0] = 0.4047206666666667
$data[3,1] = 1.214162
$data[3,2] = 0.02578034976888792
$data[3,3] = 0.02578034976888792
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.3569486666666666
$data[3,7] = 1.070846
$data[3,8] = 0.7043014606984702
$data[3,9] = 0.7043014606984702
$data[3,10] = 0.1444645023391111
$data[3,11] = 1.300180521052
$data[3,12] = 0.01815713799954523
$data[3,13] = 0.01815713799954523
$data[4,0] = 0.5058753333333333
$data[4,1] = 1.517626
$data[4,2] = 0.03222381288358415
$data[4,3] = 0.03222381288358415
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.01779166666666667
$data[4,7] = 0.053375
$data[4,8] = 0.03510503888026929
$data[4,9] = 0.03510503888026929
$data[4,10] = 0.009000365305555556
$data[4,11] = 0.08100328775
$data[4,12] = 0.001131218204148744
$data[4,13] = 0.001131218204148744
$data[5,0] = 0.5058753333333333
$data[5,1] = 1.517626
$data[5,2] = 0.03222381288358415
$data[5,3] = 0.03222381288358415
$data[5,4] = 2
$data[5,5] = 0.6666666666666666
$data[5,6] = 0.036329
$data[5,7] = 0.108987
$data[5,8] = 0.07168136529168917
$data[5,9] = 0.07168136529168917
$data[5,10] = 0.01837794498466667
$data[5,11] = 0.165401504862
$data[5,12] = 0.002309846902399235
$data[5,13] = 0.002309846902399236
$data[6,0] = 0.5058753333333333
$data[6,1] = 1.517626
$data[6,2] = 0.03222381288358415
$data[6,3] = 0.03222381288358415
$data[6,4] = 2
$data[6,5] = 0.6666666666666666
$data[6,6] = 0.09574300000000001
$data[6,7] = 0.287229
$data[6,8] = 0.1889121351295713
$data[6,9] = 0.1889121351295713
$data[6,10] = 0.04843402203933334
$data[6,11] = 0.4359061983540001
$data[6,12] = 0.00608746929385367
$data[6,13] = 0.00608746929385367
$data[7,0] = 0.5058753333333333
$data[7,1] = 1.517626
$data[7,2] = 0.03222381288358415
$data[7,3] = 0.03222381288358415
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.3569486666666666
$data[7,7] = 1.070846
$data[7,8] = 0.7043014606984702
$data[7,9] = 0.7043014606984702
$data[7,10] = 0.1805715257328889
$data[7,11] = 1.625143731596
$data[7,12] = 0.0226952784831825
$data[7,13] = 0.0226952784831825
$data[8,0] = 1.039987
$data[8,1] = 3.119961
$data[8,2] = 0.06624625531460326
$data[8,3] = 0.06624625531460326
$data[8,4] = 1
$data[8,5] = 0.3333333333333333
$data[8,6] = 0.01779166666666667
$data[8,7] = 0.053375
$data[8,8] = 0.03510503888026929
$data[8,9] = 0.03510503888026929
$data[8,10] = 0.01850310204166667
$data[8,11] = 0.166527918375
$data[8,12] = 0.002325577368491394
$data[8,13] = 0.002325577368491393
$data[9,0] = 1.039987
$data[9,1] = 3.119961
$data[9,2] = 0.06624625531460326
$data[9,3] = 0.06624625531460326
$data[9,4] = 2
$data[9,5] = 0.6666666666666666
$data[9,6] = 0.036329
$data[9,7] = 0.108987
$data[9,8] = 0.07168136529168917
$data[9,9] = 0.07168136529168917
$data[9,10] = 0.037781687723
$data[9,11] = 0.340035189507
$data[9,12] = 0.004748622026412581
$data[9,13] = 0.004748622026412581
$data[10,0] = 1.039987
$data[10,1] = 3.119961
$data[10,2] = 0.06624625531460326
$data[10,3] = 0.06624625531460326
$data[10,4] = 2
$data[10,5] = 0.6666666666666666
$data[10,6] = 0.09574300000000001
$data[10,7] = 0.287229
$data[10,8] = 0.1889121351295713
$data[10,9] = 0.1889121351295713
$data[10,10] = 0.09957147534100001
$data[10,11] = 0.8961432780690001
$data[10,12] = 0.01251472153582042
$data[10,13] = 0.01251472153582041
$data[11,0] = 1.039987
$data[11,1] = 3.119961
$data[11,2] = 0.06624625531460326
$data[11,3] = 0.06624625531460326
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 0.3569486666666666
$data[11,7] = 1.070846
$data[11,8] = 0.7043014606984702
$data[11,9] = 0.7043014606984702
$data[11,10] = 0.3712219730006666
$data[11,11] = 3.340997757006
$data[11,12] = 0.04665733438387887
$data[11,13] = 0.04665733438387887
$data[12,0] = 13.74822133333333
$data[12,1] = 41.244664
$data[12,2] = 0.8757495820329246
$data[12,3] = 0.8757495820329247
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.01779166666666667
$data[12,7] = 0.053375
$data[12,8] = 0.03510503888026929
$data[12,9] = 0.03510503888026929
$data[12,10] = 0.2446037712222222
$data[12,11] = 2.201433941
$data[12,12] = 0.0307432231266454
$data[12,13] = 0.0307432231266454
$data[13,0] = 13.74822133333333
$data[13,1] = 41.244664
$data[13,2] = 0.8757495820329246
$data[13,3] = 0.8757495820329247
$data[13,4] = 2
$data[13,5] = 0.6666666666666666
$data[13,6] = 0.036329
$data[13,7] = 0.108987
$data[13,8] = 0.07168136529168917
$data[13,9] = 0.07168136529168917
$data[13,10] = 0.4994591328186667
$data[13,11] = 4.495132195368
$data[13,12] = 0.06277492569374618
$data[13,13] = 0.06277492569374619
$data[14,0] = 13.74822133333333
$data[14,1] = 41.244664
$data[14,2] = 0.8757495820329246
$data[14,3] = 0.8757495820329247
$data[14,4] = 2
$data[14,5] = 0.6666666666666666
$data[14,6] = 0.09574300000000001
$data[14,7] = 0.287229
$data[14,8] = 0.1889121351295713
$data[14,9] = 0.1889121351295713
$data[14,10] = 1.316295955117333
$data[14,11] = 11.846663596056
$data[14,12] = 0.1654397233806695
$data[14,13] = 0.1654397233806695
$data[15,0] = 13.74822133333333
$data[15,1] = 41.244664
$data[15,2] = 0.8757495820329246
$data[15,3] = 0.8757495820329247
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 0.3569486666666666
$data[15,7] = 1.070846
$data[15,8] = 0.7043014606984702
$data[15,9] = 0.7043014606984702
$data[15,10] = 4.907409273971555
$data[15,11] = 44.166683465744
$data[15,12] = 0.6167917098318636
$data[15,13] = 0.6167917098318637

$rng = $ws.Range("G2:T17")
$rng.Value2 = $data
